$p = $ppt.ActivePresentation
$p.Slides.Item(10).Delete()
$p.Slides.Item(8).Delete()
$p.Slides.Item(7).Delete()
